$wb = $excel.ActiveWorkbook

# --- Measures sheet: insert a new column O ("item_num"), fill rows 2:25 with 1 ---
$wsMeasures = $wb.Worksheets.Item("Measures")
$wsMeasures.Columns.Item(15).Insert()
$wsMeasures.Range("O1").Value = "item_num"
$wsMeasures.Range("O2:O25").Value = 1

# --- ID sheet: insert a new column O ("item_num" header), pushing the old comment column to P ---
$wsID = $wb.Worksheets.Item("ID")
$wsID.Columns.Item(15).Insert()
$wsID.Range("O1").Value = "item_num"

# --- Dems sheet: same treatment ---
$wsDems = $wb.Worksheets.Item("Dems")
$wsDems.Columns.Item(15).Insert()
$wsDems.Range("O1").Value = "item_num"

# --- Dates sheet: same treatment ---
$wsDates = $wb.Worksheets.Item("Dates")
$wsDates.Columns.Item(15).Insert()
$wsDates.Range("O1").Value = "item_num"

# --- NewVars sheet: same treatment (header row only) ---
$wsNewVars = $wb.Worksheets.Item("NewVars")
$wsNewVars.Columns.Item(15).Insert()
$wsNewVars.Range("O1").Value = "item_num"

# --- update the _FilterDatabase defined names to cover the new column ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "ID!_FilterDatabase") {
        $n.RefersTo = "=ID!`$A`$1:`$P`$1"
    }
    if ($n.Name -eq "Measures!_FilterDatabase") {
        $n.RefersTo = "=Measures!`$A`$1:`$P`$1"
    }
}

# --- selections, matching final cursor position on each sheet ---
$wsMeasures.Range("T11").Select()
$wsID.Columns.Item(15).Select()
$wsDems.Columns.Item(15).Select()
$wsDates.Columns.Item(15).Select()
$wsNewVars.Columns.Item(15).Select()
